$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 89.8
$ws.Range("I9").Value = 99.666664
$ws.Range("K9").Value = 99.666664
$ws.Range("M9").Value = 69.333336

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 759.8889
$ws.Range("I33").Value = 408.33334
$ws.Range("J33").Value = 1463
$ws.Range("K33").Value = 408.33334
$ws.Range("L33").Value = 1463
$ws.Range("M33").Value = -179.33334
$ws.Range("N33").Value = -1921

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5263.143
$ws.Range("J62").Value = 5716.5
$ws.Range("L62").Value = 5716.5
$ws.Range("N62").Value = -6964.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 5263.143
$ws.Range("J65").Value = 5716.5
$ws.Range("L65").Value = 28582.5
$ws.Range("N65").Value = -34822.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 128732.52
$ws.Range("J69").Value = 139578.88
$ws.Range("L69").Value = 418736.64
$ws.Range("N69").Value = -420484.64

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 128732.52
$ws.Range("J72").Value = 139578.88
$ws.Range("L72").Value = 1256209.92
$ws.Range("N72").Value = -1264945.92

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2669.3635
$ws.Range("I106").Value = 2066.3
$ws.Range("K106").Value = 2066.3
$ws.Range("M106").Value = -1435.3

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5496584.5
$ws.Range("I132").Value = 6804467
$ws.Range("K132").Value = 20413401
$ws.Range("M132").Value = -20410871

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 16646.953
$ws.Range("I137").Value = 1954.75
$ws.Range("K137").Value = 5864.25
$ws.Range("M137").Value = -3314.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1660781.9
$ws.Range("I138").Value = 2403
$ws.Range("J138").Value = 2690120.5
$ws.Range("K138").Value = 7209
$ws.Range("L138").Value = 8070361.5
$ws.Range("M138").Value = -2069
$ws.Range("N138").Value = -8080641.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 445.4
$ws.Range("I2").Value = 383.66666
$ws.Range("K2").Value = 383.66666
$ws.Range("M2").Value = -270.66666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17552.045
$ws.Range("I32").Value = 17973.047
$ws.Range("K32").Value = 17973.047
$ws.Range("M32").Value = -17686.047

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4891.5835
$ws.Range("I45").Value = 3491.8333
$ws.Range("J45").Value = 6291.3335
$ws.Range("K45").Value = 3491.8333
$ws.Range("L45").Value = 6291.3335
$ws.Range("M45").Value = -3114.8333
$ws.Range("N45").Value = -7045.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5808.878
$ws.Range("I61").Value = 3381.8286
$ws.Range("K61").Value = 3381.8286
$ws.Range("M61").Value = -3169.8286

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3092.9648
$ws.Range("I74").Value = 1122.7441
$ws.Range("K74").Value = 1122.7441
$ws.Range("M74").Value = -248.7440999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3092.9648
$ws.Range("I77").Value = 1122.7441
$ws.Range("K77").Value = 5613.720499999999
$ws.Range("M77").Value = -1245.720499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 445.4
$ws.Range("I116").Value = 383.66666
$ws.Range("K116").Value = 383.66666
$ws.Range("M116").Value = 1910.33334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5808.878
$ws.Range("I136").Value = 3381.8286
$ws.Range("K136").Value = 10145.4858
$ws.Range("M136").Value = -7595.485799999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 445.4
$ws.Range("I3").Value = 383.66666
$ws.Range("K3").Value = 383.66666
$ws.Range("M3").Value = -269.66666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 50000
$ws.Range("J56").Value = 50000
$ws.Range("L56").Value = 50000
$ws.Range("N56").Value = -51478

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2793.7273
$ws.Range("I134").Value = 2685.65
$ws.Range("J134").Value = 3874.5
$ws.Range("K134").Value = 8056.950000000001
$ws.Range("L134").Value = 11623.5
$ws.Range("M134").Value = -5521.950000000001
$ws.Range("N134").Value = -16693.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 14926.167
$ws.Range("I6").Value = 21389.25
$ws.Range("K6").Value = 21389.25
$ws.Range("M6").Value = -21276.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5887498.5
$ws.Range("I31").Value = 14291512
$ws.Range("J31").Value = 4689.1
$ws.Range("K31").Value = 14291512
$ws.Range("L31").Value = 4689.1
$ws.Range("M31").Value = -14291217
$ws.Range("N31").Value = -5279.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5887498.5
$ws.Range("I34").Value = 14291512
$ws.Range("J34").Value = 4689.1
$ws.Range("K34").Value = 14291512
$ws.Range("L34").Value = 4689.1
$ws.Range("M34").Value = -14291310
$ws.Range("N34").Value = -5093.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 9382.444
$ws.Range("I103").Value = 9382.444
$ws.Range("K103").Value = 9382.444
$ws.Range("M103").Value = -8210.444

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1817.9286
$ws.Range("I105").Value = 1731.125
$ws.Range("K105").Value = 1731.125
$ws.Range("M105").Value = 15.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 408.55554
$ws.Range("I107").Value = 248.95238
$ws.Range("K107").Value = 248.95238
$ws.Range("M107").Value = 1671.04762

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1668.7142
$ws.Range("I132").Value = 1461.7916
$ws.Range("J132").Value = 2910.25
$ws.Range("K132").Value = 4385.3748
$ws.Range("L132").Value = 8730.75
$ws.Range("M132").Value = -1855.3748
$ws.Range("N132").Value = -13790.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 119805.63
$ws.Range("J140").Value = 119805.63
$ws.Range("L140").Value = 119805.63
$ws.Range("N140").Value = -130165.63

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 595.1739
$ws.Range("I5").Value = 535.6667
$ws.Range("J5").Value = 633.4286
$ws.Range("K5").Value = 1607.0001
$ws.Range("L5").Value = 1900.2858
$ws.Range("M5").Value = -1495.0001
$ws.Range("N5").Value = -2124.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6277.364
$ws.Range("I56").Value = 6277.364
$ws.Range("K56").Value = 6277.364
$ws.Range("M56").Value = -5747.364

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4071.1035
$ws.Range("J68").Value = 4423.32
$ws.Range("L68").Value = 13269.96
$ws.Range("N68").Value = -14891.96

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 4071.1035
$ws.Range("J71").Value = 4423.32
$ws.Range("L71").Value = 39809.88
$ws.Range("N71").Value = -47921.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 4420.4
$ws.Range("J129").Value = 4093
$ws.Range("L129").Value = 12279
$ws.Range("N129").Value = -22279

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1355.45
$ws.Range("I132").Value = 1227.875
$ws.Range("K132").Value = 11050.875
$ws.Range("M132").Value = -8520.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 595.1739
$ws.Range("I135").Value = 535.6667
$ws.Range("J135").Value = 633.4286
$ws.Range("K135").Value = 4821.0003
$ws.Range("L135").Value = 5700.8574
$ws.Range("M135").Value = -2286.0003
$ws.Range("N135").Value = -10770.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2376.25
$ws.Range("I102").Value = 2001.8334
$ws.Range("J102").Value = 3499.5
$ws.Range("K102").Value = 2001.8334
$ws.Range("L102").Value = 3499.5
$ws.Range("M102").Value = -379.8334
$ws.Range("N102").Value = -6743.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7822.56
$ws.Range("I132").Value = 8030
$ws.Range("K132").Value = 24090
$ws.Range("M132").Value = -21560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2886.9092
$ws.Range("I40").Value = 2782.8386
$ws.Range("K40").Value = 2782.8386
$ws.Range("M40").Value = -2646.8386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4514.533
$ws.Range("I132").Value = 3837.3845
$ws.Range("J132").Value = 5032.353
$ws.Range("K132").Value = 11512.1535
$ws.Range("L132").Value = 15097.059
$ws.Range("M132").Value = -8982.1535
$ws.Range("N132").Value = -20157.059

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4859.353
$ws.Range("I136").Value = 2599
$ws.Range("J136").Value = 6868.5557
$ws.Range("K136").Value = 7797
$ws.Range("L136").Value = 20605.6671
$ws.Range("M136").Value = -5247
$ws.Range("N136").Value = -25705.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5173.483
$ws.Range("I81").Value = 4940.6665
$ws.Range("K81").Value = 9881.333000000001
$ws.Range("M81").Value = -8820.333000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 5173.483
$ws.Range("I84").Value = 4940.6665
$ws.Range("K84").Value = 49406.665
$ws.Range("M84").Value = -44102.665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 13396.333
$ws.Range("I88").Value = 6000
$ws.Range("J88").Value = 17094.5
$ws.Range("K88").Value = 6000
$ws.Range("L88").Value = 17094.5
$ws.Range("M88").Value = -5594
$ws.Range("N88").Value = -17906.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H91").Value = 13396.333
$ws.Range("I91").Value = 6000
$ws.Range("J91").Value = 17094.5
$ws.Range("K91").Value = 6000
$ws.Range("L91").Value = 17094.5
$ws.Range("M91").Value = -4596
$ws.Range("N91").Value = -19902.5
